$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (status "In-progress") keeps its text but gets a percentage number format applied.
$ws.Range("G3").NumberFormat = "0%"

# Row 4 status changes from text "In-progress" to a numeric value of 25%.
$ws.Range("G4").NumberFormat = "0%"
$ws.Range("G4").Value = 0.25

# Update the active cell selection to D7.
$ws.Range("D7").Select()
